# New crime data collected - weekly CompStat report update (63rd Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Header text: volume/issue number and the reporting week date range
# -------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# -------------------------------------------------------------------
# A handful of data cells flip between a numeric value and the
# special "0" / "***.*" placeholder text used elsewhere in the sheet
# (e.g. row 22, which is entirely placeholders). Setting .Value alone
# would just coerce "0" back into a literal 0, and it would stamp a
# brand-new number format instead of reusing the existing text style,
# so instead we Copy a same-styled donor cell (which duplicates both
# format and content) and then overwrite the value where it must
# differ from the donor.
#   C22 -> style "14", shared text "0"
#   E22 -> style "14", shared text "***.*"
#   F15 -> style "15" (plain integer number format)
#   H15 -> style "16" (one-decimal number format)
# -------------------------------------------------------------------

# Row 15: Murder - C,D -> "0", E -> "***.*"
$ws.Range("C22").Copy($ws.Range("C15"))
$ws.Range("C22").Copy($ws.Range("D15"))
$ws.Range("E22").Copy($ws.Range("E15"))

# Row 16: Rape - D -> "0", E -> "***.*"
$ws.Range("C22").Copy($ws.Range("D16"))
$ws.Range("E22").Copy($ws.Range("E16"))

# Row 18: Burglary - C -> "0"; D,E -> numeric
$ws.Range("C22").Copy($ws.Range("C18"))
$ws.Range("F15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 2
$ws.Range("H15").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -100

# Row 23: Transit - C -> "0"; D,E -> numeric
$ws.Range("C22").Copy($ws.Range("C23"))
$ws.Range("F15").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("H15").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100

# Row 26: G.L.A. - C,D -> "0", E -> "***.*"
$ws.Range("C22").Copy($ws.Range("C26"))
$ws.Range("C22").Copy($ws.Range("D26"))
$ws.Range("E22").Copy($ws.Range("E26"))

# -------------------------------------------------------------------
# Remaining cells: plain numeric value updates (counts and the
# derived percent-change figures recomputed from the new counts).
# -------------------------------------------------------------------

# Row 15: Murder
$ws.Range("N15").Value = -68.75

# Row 16: Rape
$ws.Range("C16").Value = 6
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = 120
$ws.Range("I16").Value = 54
$ws.Range("K16").Value = 1.886792452830
$ws.Range("L16").Value = 58.823529411764
$ws.Range("M16").Value = -53.846153846153
$ws.Range("N16").Value = -84.302325581395

# Row 17: Robbery
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 66.666666666666
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 88
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = 6.024096385542
$ws.Range("L17").Value = 17.333333333333
$ws.Range("M17").Value = 41.935483870967
$ws.Range("N17").Value = -44.654088050314

# Row 18: Fel. Assault
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 49
$ws.Range("K18").Value = 2.040816326530
$ws.Range("L18").Value = -7.407407407407
$ws.Range("M18").Value = -59.016393442622
$ws.Range("N18").Value = -91.909385113268

# Row 19: Burglary
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 85.714285714285
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 299
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = 24.583333333333
$ws.Range("L19").Value = 86.875
$ws.Range("M19").Value = 25.630252100840
$ws.Range("N19").Value = -11.538461538461

# Row 20: Gr. Larceny
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 37.5
$ws.Range("I20").Value = 55
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -3.508771929824
$ws.Range("L20").Value = 30.952380952381
$ws.Range("M20").Value = -33.734939759036
$ws.Range("N20").Value = -96.099290780141

# Row 21: TOTAL
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 86.666666666666
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 20.512820512820
$ws.Range("I21").Value = 553
$ws.Range("J21").Value = 488
$ws.Range("K21").Value = 13.319672131147
$ws.Range("L21").Value = 49.056603773584
$ws.Range("M21").Value = -12.222222222222
$ws.Range("N21").Value = -80.898100172711

# Row 23: Transit
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -83.333333333333
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = -60
$ws.Range("L23").Value = 14.285714285714
$ws.Range("M23").Value = -52.941176470588

# Row 24: Housing
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 37.5
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 4.255319148936
$ws.Range("I24").Value = 598
$ws.Range("J24").Value = 536
$ws.Range("K24").Value = 11.567164179104
$ws.Range("L24").Value = 56.544502617801
$ws.Range("M24").Value = 33.184855233853

# Row 25: Petit Larceny
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 5.263157894736
$ws.Range("I25").Value = 133
$ws.Range("J25").Value = 109
$ws.Range("K25").Value = 22.018348623853
$ws.Range("L25").Value = 12.711864406779
$ws.Range("M25").Value = -22.674418604651

# Row 26: Misd. Assault
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = 30

# Row 27: UCR Rape*
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -5.882352941176

# Row 28: Other Sex Crimes
$ws.Range("N28").Value = -81.481481481481

# Row 29: Shooting Vic.
$ws.Range("N29").Value = -83.333333333333

# Row 30: Shooting Inc.
$ws.Range("L30").Value = -50
